$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Scores" label and AVERAGE formulas in row 32
$ws.Range("A32").Value = "Scores"
$ws.Range("B32").Formula = "=AVERAGE(B2:B31)"
$ws.Range("C32").Formula = "=AVERAGE(C2:C31)"
$ws.Range("D32").Formula = "=AVERAGE(D2:D31)"
$ws.Range("E32").Formula = "=AVERAGE(E2:E31)"
$ws.Range("F32").Formula = "=AVERAGE(F2:F31)"
$ws.Range("G32").Formula = "=AVERAGE(G2:G31)"

# Apply border style to A1:G32 (new style with border only)
$ws.Range("A1:G32").Borders.LineStyle = 1

# Select cell I12 to match the saved view state
$ws.Range("I12").Select() | Out-Null

# Match the saved print/page setup (paper size + portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
